$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the URL value (matchsource -> matchsync)
$ws.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-disease-codes"

# Set the Experimental value (was blank) to the literal text "true".
# A direct Value/Formula assignment of the bare word "true" is auto-coerced
# to a native Boolean by the engine (like typing it into Excel), so instead
# build it as a formula that evaluates to the text string "true" and then
# bake the formula result down to a plain value with a values-only paste.
# This keeps the cell's existing style and produces a genuine text cell.
$b7 = $ws.Range("B7")
$b7.Formula = '=TEXT("true","@")'
$b7.Copy()
$b7.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# Update the Date value to the new timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
